# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Change: cell B11 on the "Rules" sheet used to hold the text "R40"
# (rule id). It now holds the text "1" instead, while keeping the
# cell's existing look (style/border/etc.) exactly as-is.
#
# A plain  $ws.Range("B11").Value = "1"  would make Excel's normal
# "smart" type detection store it as the *number* 1 (General format),
# which also drops the cell out of the shared-string table. To keep it
# a genuine text value (so it round-trips as t="s" in the XML, matching
# the original text-cell typing) without disturbing B11's style, we
# stage the text in a scratch cell that is explicitly formatted as
# Text, copy it, and paste only the *value* (xlPasteValues) into B11 -
# that leaves B11's own formatting/style untouched and only replaces
# its stored value/type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "1"

$scratch.Copy()
$ws.Range("B11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$scratch.Clear()
$excel.CutCopyMode = 0
